$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 20. This shifts the former row 20 -> 21 and
# former row 21 -> 22 (both keep their original values/formatting intact).
$ws.Rows(20).Insert()

# Populate the newly inserted row 20 with this week's new data entry.
$ws.Range("A20").Value = 7
$ws.Range("B20").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C20").Value = "Ñuble"
$ws.Range("D20").Value = 44782
$ws.Range("E20").Value = 16
$ws.Range("F20").Value = 100112043
$ws.Range("G20").Value = "Pepino dulce"
$ws.Range("H20").Value = "Cultivar IV Región"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 120
$ws.Range("K20").Value = 17000
$ws.Range("L20").Value = 18000
$ws.Range("M20").Value = 17500
$ws.Range("N20").Value = "$/bandeja 18 kilos"
$ws.Range("O20").Value = "Provincia de Limarí"
$ws.Range("P20").Value = 972
$ws.Range("Q20").Value = 18
$ws.Range("R20").Value = "Hortaliza"
